# RSD Retrofit: Use LimType FX with NCAP_AF
# Insert a new "LimType" column (value "FX") into the rsd_rtft table on the
# RSD_RTFT sheet, between "Attribute" and "Year", and make RSD_RTFT the
# active sheet/tab with the selection on J22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSD_RTFT")

# Shift columns D:G (Year, Pset_PN, IE, National) one column to the right,
# opening up column D for the new LimType column.
$ws.Range("D1").EntireColumn.Insert()

# Populate the new LimType column (header + the one data row).
$ws.Range("D2").Value = "LimType"
$ws.Range("D3").Value = "FX"

# Re-assert the header labels that shifted right, so the engine treats every
# header cell as freshly written (keeps the table's column-name sync honest).
$ws.Range("E2").Value = "Year"
$ws.Range("F2").Value = "Pset_PN"
$ws.Range("G2").Value = "IE"
$ws.Range("H2").Value = "National"

# Grow the table (ListObject) to cover the new column.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B2:H3"))

# The rightmost column in the resized table keeps a synthetic default name
# ("ColumnN") instead of picking up the "National" header text, so fix it up
# by writing straight to its header cell through the table's header range.
$tbl.HeaderRowRange.Cells.Item(1, 7).Value = "National"

# Make RSD_RTFT the active sheet/tab with the same selection as the source
# edit (activeCell J22).
$ws.Activate()
$ws.Range("J22").Select()
